$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1677.7778
$ws.Range("H86").Value = 50727.25
$ws.Range("I86").Value = 91954.5
$ws.Range("J86").Value = 9500
$ws.Range("K86").Value = 91954.5
$ws.Range("L86").Value = 9500
$ws.Range("M86").Value = -90831.5
$ws.Range("N86").Value = -11746
$ws.Range("H87").Value = 39985
$ws.Range("J87").Value = 39985
$ws.Range("L87").Value = 39985
$ws.Range("N87").Value = -42481
$ws.Range("H89").Value = 50727.25
$ws.Range("I89").Value = 91954.5
$ws.Range("J89").Value = 9500
$ws.Range("K89").Value = 459772.5
$ws.Range("L89").Value = 47500
$ws.Range("M89").Value = -454156.5
$ws.Range("N89").Value = -58732
$ws.Range("H90").Value = 39985
$ws.Range("J90").Value = 39985
$ws.Range("L90").Value = 119955
$ws.Range("N90").Value = -132435
$ws.Range("H135").Value = 821
$ws.Range("I135").Value = 795.4
$ws.Range("K135").Value = 7158.599999999999
$ws.Range("M135").Value = -4623.599999999999
$ws.Range("H138").Value = 2653.2173
$ws.Range("J138").Value = 4076.923
$ws.Range("L138").Value = 12230.769
$ws.Range("N138").Value = -22510.769

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1666.9
$ws.Range("I2").Value = 770.375
$ws.Range("J2").Value = 5253
$ws.Range("K2").Value = 770.375
$ws.Range("L2").Value = 5253
$ws.Range("M2").Value = -657.375
$ws.Range("N2").Value = -5479
$ws.Range("H45").Value = 2571.5
$ws.Range("I45").Value = 2577
$ws.Range("K45").Value = 2577
$ws.Range("M45").Value = -2200
$ws.Range("H55").Value = 20000
$ws.Range("J55").Value = 20000
$ws.Range("L55").Value = 20000
$ws.Range("N55").Value = -20630
$ws.Range("H61").Value = 1769.909
$ws.Range("I61").Value = 1646.9
$ws.Range("K61").Value = 1646.9
$ws.Range("M61").Value = -1434.9
$ws.Range("H88").Value = 2172.5386
$ws.Range("I88").Value = 1189.3334
$ws.Range("J88").Value = 3015.2856
$ws.Range("K88").Value = 1189.3334
$ws.Range("L88").Value = 3015.2856
$ws.Range("M88").Value = -783.3334
$ws.Range("N88").Value = -3827.2856
$ws.Range("H91").Value = 2172.5386
$ws.Range("I91").Value = 1189.3334
$ws.Range("J91").Value = 3015.2856
$ws.Range("K91").Value = 1189.3334
$ws.Range("L91").Value = 3015.2856
$ws.Range("M91").Value = 214.6666
$ws.Range("N91").Value = -5823.2856
$ws.Range("H116").Value = 1666.9
$ws.Range("I116").Value = 770.375
$ws.Range("J116").Value = 5253
$ws.Range("K116").Value = 770.375
$ws.Range("L116").Value = 5253
$ws.Range("M116").Value = 1523.625
$ws.Range("N116").Value = -9841
$ws.Range("H132").Value = 2491.8
$ws.Range("I132").Value = 2491.8
$ws.Range("K132").Value = 7475.400000000001
$ws.Range("M132").Value = -4945.400000000001
$ws.Range("H136").Value = 1769.909
$ws.Range("I136").Value = 1646.9
$ws.Range("K136").Value = 4940.700000000001
$ws.Range("M136").Value = -2390.700000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1666.9
$ws.Range("I3").Value = 770.375
$ws.Range("J3").Value = 5253
$ws.Range("K3").Value = 770.375
$ws.Range("L3").Value = 5253
$ws.Range("M3").Value = -656.375
$ws.Range("N3").Value = -5481
$ws.Range("H82").Value = 21047.25
$ws.Range("J82").Value = 30000
$ws.Range("L82").Value = 30000
$ws.Range("N82").Value = -30766
$ws.Range("H85").Value = 21047.25
$ws.Range("J85").Value = 30000
$ws.Range("L85").Value = 30000
$ws.Range("N85").Value = -32652
$ws.Range("H134").Value = 6918.4736
$ws.Range("I134").Value = 6275.0557
$ws.Range("K134").Value = 18825.1671
$ws.Range("M134").Value = -16290.1671

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1700.875
$ws.Range("I31").Value = 1498.75
$ws.Range("K31").Value = 1498.75
$ws.Range("M31").Value = -1203.75
$ws.Range("H34").Value = 1700.875
$ws.Range("I34").Value = 1498.75
$ws.Range("K34").Value = 1498.75
$ws.Range("M34").Value = -1296.75
$ws.Range("H58").Value = 1606.8462
$ws.Range("I58").Value = 1631
$ws.Range("J58").Value = 1589.1333
$ws.Range("K58").Value = 1631
$ws.Range("L58").Value = 1589.1333
$ws.Range("M58").Value = -1428
$ws.Range("N58").Value = -1995.1333
$ws.Range("H59").Value = 29084
$ws.Range("J59").Value = 30000
$ws.Range("L59").Value = 30000
$ws.Range("N59").Value = -32290
$ws.Range("H74").Value = 29997.727
$ws.Range("J74").Value = 29997.727
$ws.Range("L74").Value = 29997.727
$ws.Range("N74").Value = -31745.727
$ws.Range("H77").Value = 29997.727
$ws.Range("J77").Value = 29997.727
$ws.Range("L77").Value = 89993.181
$ws.Range("N77").Value = -98729.181
$ws.Range("H94").Value = 1011.25
$ws.Range("I94").Value = 1011.25
$ws.Range("K94").Value = 1011.25
$ws.Range("M94").Value = -560.25
$ws.Range("H105").Value = 687.6667
$ws.Range("I105").Value = 698.3333
$ws.Range("K105").Value = 698.3333
$ws.Range("M105").Value = 1048.6667
$ws.Range("H122").Value = 833.3333
$ws.Range("I122").Value = 833.3333
$ws.Range("K122").Value = 2499.9999
$ws.Range("M122").Value = -49.9998999999998
$ws.Range("H132").Value = 3404.111
$ws.Range("I132").Value = 3457.125
$ws.Range("J132").Value = 2980
$ws.Range("K132").Value = 10371.375
$ws.Range("L132").Value = 8940
$ws.Range("M132").Value = -7841.375
$ws.Range("N132").Value = -14000
$ws.Range("H134").Value = 4424.304
$ws.Range("I134").Value = 4126.25
$ws.Range("K134").Value = 12378.75
$ws.Range("M134").Value = -9843.75
$ws.Range("H136").Value = 1606.8462
$ws.Range("I136").Value = 1631
$ws.Range("J136").Value = 1589.1333
$ws.Range("K136").Value = 4893
$ws.Range("L136").Value = 4767.3999
$ws.Range("M136").Value = -2343
$ws.Range("N136").Value = -9867.3999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 990.6667
$ws.Range("I97").Value = 990.6667
$ws.Range("K97").Value = 990.6667
$ws.Range("M97").Value = -494.6667
$ws.Range("H107").Value = 391.16666
$ws.Range("I107").Value = 349.6
$ws.Range("K107").Value = 349.6
$ws.Range("M107").Value = 1570.4
$ws.Range("H123").Value = 50001
$ws.Range("J123").Value = 50001
$ws.Range("L123").Value = 50001
$ws.Range("N123").Value = -54901
$ws.Range("H132").Value = 5921
$ws.Range("I132").Value = 2894.6667
$ws.Range("K132").Value = 8684.000100000001
$ws.Range("M132").Value = -6154.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4641.1377
$ws.Range("J46").Value = 2704.2222
$ws.Range("L46").Value = 2704.2222
$ws.Range("N46").Value = -3080.2222
$ws.Range("H82").Value = 2000
$ws.Range("I82").Value = 2000
$ws.Range("K82").Value = 2000
$ws.Range("M82").Value = -1639
$ws.Range("H85").Value = 2000
$ws.Range("I85").Value = 2000
$ws.Range("K85").Value = 2000
$ws.Range("M85").Value = -752

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1783.75
$ws.Range("I96").Value = 970
$ws.Range("J96").Value = 2597.5
$ws.Range("K96").Value = 970
$ws.Range("L96").Value = 2597.5
$ws.Range("M96").Value = 403
$ws.Range("N96").Value = -5343.5
$ws.Range("H132").Value = 1999
$ws.Range("I132").Value = 1999
$ws.Range("K132").Value = 5997
$ws.Range("M132").Value = -3467
